$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$caseQuery = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.ethnicity =`"NOT_HISPANIC_OR_LATINO`" `nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

$statQuery = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.ethnicity =`"NOT_HISPANIC_OR_LATINO`" `nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

$ws.Range("C2").Value = $statQuery
$ws.Range("B2").Value = $caseQuery

$ws.Range("B7").Select()
